$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 6
$ws.Range("E6").Value = 32

# Row 25
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 76
$ws.Range("H25").Value = 76

# Row 30
$ws.Range("E30").Value = 139

# Row 41
$ws.Range("F41").Value = 114
$ws.Range("H41").Value = 114

# Row 44
$ws.Range("F44").Value = 107
$ws.Range("H44").Value = 107

# Row 45
$ws.Range("E45").Value = 75

# Row 46
$ws.Range("E46").Value = 191

# Row 47
$ws.Range("E47").Value = 289

# Row 48
$ws.Range("E48").Value = 133

# Row 49
$ws.Range("E49").Value = 166
